$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 204.3437653333333
$ws.Range("H2").Value = 613.0312959999999
$ws.Range("I2").Value = 0.1758767975891587
$ws.Range("J2").Value = 0.1775259368507247
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 172.9098156666666
$ws.Range("N2").Value = 518.7294469999999
$ws.Range("O2").Value = 0.4524900325013766
$ws.Range("P2").Value = 0.4688259813160829
$ws.Range("Q2").Value = 35333.04279641925
$ws.Range("R2").Value = 317997.3851677732
$ws.Range("S2").Value = 0.07958249785735647
$ws.Range("T2").Value = 0.08322877155309796
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 204.3437653333333
$ws.Range("H3").Value = 613.0312959999999
$ws.Range("I3").Value = 0.1758767975891587
$ws.Range("J3").Value = 0.1775259368507247
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 46.15041300000001
$ws.Range("N3").Value = 138.451239
$ws.Range("O3").Value = 0.120771639237527
$ws.Range("P3").Value = 0.1251317779701883
$ws.Range("Q3").Value = 9430.549164108415
$ws.Range("R3").Value = 84874.94247697573
$ws.Range("S3").Value = 0.02124092914868943
$ws.Range("T3").Value = 0.02221413611395455
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 204.3437653333333
$ws.Range("H4").Value = 613.0312959999999
$ws.Range("I4").Value = 0.1758767975891587
$ws.Range("J4").Value = 0.1775259368507247
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 64.30983099999999
$ws.Range("N4").Value = 192.929493
$ws.Range("O4").Value = 0.1682932655219863
$ws.Range("P4").Value = 0.1743690461446646
$ws.Range("Q4").Value = 13141.31301449032
$ws.Range("R4").Value = 118271.8171304129
$ws.Range("S4").Value = 0.02959888059582892
$ws.Range("T4").Value = 0.03095502827459883
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 204.3437653333333
$ws.Range("H5").Value = 613.0312959999999
$ws.Range("I5").Value = 0.1758767975891587
$ws.Range("J5").Value = 0.1775259368507247
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 58.81429933333334
$ws.Range("N5").Value = 176.442898
$ws.Range("O5").Value = 0.1539119344629323
$ws.Range("P5").Value = 0.1594685153879524
$ws.Range("Q5").Value = 12018.33538121509
$ws.Range("R5").Value = 108165.0184309358
$ws.Range("S5").Value = 0.027069538144093
$ws.Range("T5").Value = 0.02830979759244046
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 204.3437653333333
$ws.Range("H6").Value = 613.0312959999999
$ws.Range("I6").Value = 0.1758767975891587
$ws.Range("J6").Value = 0.1775259368507247
$ws.Range("K6").Value = 2
$ws.Range("M6").Value = 39.945198
$ws.Range("N6").Value = 79.890396
$ws.Range("O6").Value = 0.1045331282761778
$ws.Range("P6").Value = 0.07220467918111168
$ws.Range("Q6").Value = 8162.552166305534
$ws.Range("R6").Value = 48975.31299783321
$ws.Range("S6").Value = 0.01838495184319089
$ws.Range("T6").Value = 0.01281820331663287
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 161.8155033333333
$ws.Range("H7").Value = 485.44651
$ws.Range("I7").Value = 0.1392731140102079
$ws.Range("J7").Value = 0.1405790324914581
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 172.9098156666666
$ws.Range("N7").Value = 518.7294469999999
$ws.Range("O7").Value = 0.4524900325013766
$ws.Range("P7").Value = 0.4688259813160829
$ws.Range("Q7").Value = 27979.48885337555
$ws.Range("R7").Value = 251815.3996803799
$ws.Range("S7").Value = 0.06301969588504691
$ws.Range("T7").Value = 0.06590710286027338
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 161.8155033333333
$ws.Range("H8").Value = 485.44651
$ws.Range("I8").Value = 0.1392731140102079
$ws.Range("J8").Value = 0.1405790324914581
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 46.15041300000001
$ws.Range("N8").Value = 138.451239
$ws.Range("O8").Value = 0.120771639237527
$ws.Range("P8").Value = 0.1251317779701883
$ws.Range("Q8").Value = 7467.852308636211
$ws.Range("R8").Value = 67210.6707777259
$ws.Range("S8").Value = 0.0168202422807278
$ws.Range("T8").Value = 0.01759090428098503
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 161.8155033333333
$ws.Range("H9").Value = 485.44651
$ws.Range("I9").Value = 0.1392731140102079
$ws.Range("J9").Value = 0.1405790324914581
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 64.30983099999999
$ws.Range("N9").Value = 192.929493
$ws.Range("O9").Value = 0.1682932655219863
$ws.Range("P9").Value = 0.1743690461446646
$ws.Range("Q9").Value = 10406.3276725466
$ws.Range("R9").Value = 93656.94905291942
$ws.Range("S9").Value = 0.02343872715619379
$ws.Range("T9").Value = 0.02451263180347537
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 161.8155033333333
$ws.Range("H10").Value = 485.44651
$ws.Range("I10").Value = 0.1392731140102079
$ws.Range("J10").Value = 0.1405790324914581
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 58.81429933333334
$ws.Range("N10").Value = 176.442898
$ws.Range("O10").Value = 0.1539119344629323
$ws.Range("P10").Value = 0.1594685153879524
$ws.Range("Q10").Value = 9517.065449820666
$ws.Range("R10").Value = 85653.58904838598
$ws.Range("S10").Value = 0.02143579439598762
$ws.Range("T10").Value = 0.02241792960608755
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 161.8155033333333
$ws.Range("H11").Value = 485.44651
$ws.Range("I11").Value = 0.1392731140102079
$ws.Range("J11").Value = 0.1405790324914581
$ws.Range("K11").Value = 2
$ws.Range("M11").Value = 39.945198
$ws.Range("N11").Value = 79.890396
$ws.Range("O11").Value = 0.1045331282761778
$ws.Range("P11").Value = 0.07220467918111168
$ws.Range("Q11").Value = 6463.75232011966
$ws.Range("R11").Value = 38782.51392071796
$ws.Range("S11").Value = 0.0145586542922518
$ws.Range("T11").Value = 0.01015046394063681
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 467.064364
$ws.Range("H12").Value = 1401.193092
$ws.Range("I12").Value = 0.4019979981984663
$ws.Range("J12").Value = 0.4057674020708785
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 172.9098156666666
$ws.Range("N12").Value = 518.7294469999999
$ws.Range("O12").Value = 0.4524900325013766
$ws.Range("P12").Value = 0.4688259813160829
$ws.Range("Q12").Value = 80760.01308370889
$ws.Range("R12").Value = 726840.11775338
$ws.Range("S12").Value = 0.1819000872703124
$ws.Range("T12").Value = 0.1902343004619572
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 467.064364
$ws.Range("H13").Value = 1401.193092
$ws.Range("I13").Value = 0.4019979981984663
$ws.Range("J13").Value = 0.4057674020708785
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 46.15041300000001
$ws.Range("N13").Value = 138.451239
$ws.Range("O13").Value = 0.120771639237527
$ws.Range("P13").Value = 0.1251317779701883
$ws.Range("Q13").Value = 21555.21329618234
$ws.Range("R13").Value = 193996.919665641
$ws.Range("S13").Value = 0.0485499572126332
$ws.Range("T13").Value = 0.05077439646347329
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 467.064364
$ws.Range("H14").Value = 1401.193092
$ws.Range("I14").Value = 0.4019979981984663
$ws.Range("J14").Value = 0.4057674020708785
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 64.30983099999999
$ws.Range("N14").Value = 192.929493
$ws.Range("O14").Value = 0.1682932655219863
$ws.Range("P14").Value = 0.1743690461446646
$ws.Range("Q14").Value = 30036.83031496248
$ws.Range("R14").Value = 270331.4728346623
$ws.Range("S14").Value = 0.06765355585012146
$ws.Range("T14").Value = 0.07075327485569768
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 467.064364
$ws.Range("H15").Value = 1401.193092
$ws.Range("I15").Value = 0.4019979981984663
$ws.Range("J15").Value = 0.4057674020708785
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 58.81429933333334
$ws.Range("N15").Value = 176.442898
$ws.Range("O15").Value = 0.1539119344629323
$ws.Range("P15").Value = 0.1594685153879524
$ws.Range("Q15").Value = 27470.06331222896
$ws.Range("R15").Value = 247230.5698100606
$ws.Range("S15").Value = 0.06187228955295231
$ws.Range("T15").Value = 0.06470712520106935
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 467.064364
$ws.Range("H16").Value = 1401.193092
$ws.Range("I16").Value = 0.4019979981984663
$ws.Range("J16").Value = 0.4057674020708785
$ws.Range("K16").Value = 2
$ws.Range("M16").Value = 39.945198
$ws.Range("N16").Value = 79.890396
$ws.Range("O16").Value = 0.1045331282761778
$ws.Range("P16").Value = 0.07220467918111168
$ws.Range("Q16").Value = 18656.97849872407
$ws.Range("R16").Value = 111941.8709923444
$ws.Range("S16").Value = 0.04202210831244697
$ws.Range("T16").Value = 0.02929830508868093
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 296.2543436666667
$ws.Range("H17").Value = 888.763031
$ws.Range("I17").Value = 0.2549833862118423
$ws.Range("J17").Value = 0.2573742821060879
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 172.9098156666666
$ws.Range("N17").Value = 518.7294469999999
$ws.Range("O17").Value = 0.4524900325013766
$ws.Range("P17").Value = 0.4688259813160829
$ws.Range("Q17").Value = 51225.28395385265
$ws.Range("R17").Value = 461027.5555846738
$ws.Range("S17").Value = 0.1153774407143076
$ws.Range("T17").Value = 0.120663750373909
$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 296.2543436666667
$ws.Range("H18").Value = 888.763031
$ws.Range("I18").Value = 0.2549833862118423
$ws.Range("J18").Value = 0.2573742821060879
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 46.15041300000001
$ws.Range("N18").Value = 138.451239
$ws.Range("O18").Value = 0.120771639237527
$ws.Range("P18").Value = 0.1251317779701883
$ws.Range("Q18").Value = 13672.2603132606
$ws.Range("R18").Value = 123050.3428193454
$ws.Range("S18").Value = 0.03079476153113963
$ws.Range("T18").Value = 0.0322057015237356
$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 296.2543436666667
$ws.Range("H19").Value = 888.763031
$ws.Range("I19").Value = 0.2549833862118423
$ws.Range("J19").Value = 0.2573742821060879
$ws.Range("K19").Value = 3
$ws.Range("M19").Value = 64.30983099999999
$ws.Range("N19").Value = 192.929493
$ws.Range("O19").Value = 0.1682932655219863
$ws.Range("P19").Value = 0.1743690461446646
$ws.Range("Q19").Value = 19052.06677421925
$ws.Range("R19").Value = 171468.6009679733
$ws.Range("S19").Value = 0.04291198671944475
$ws.Range("T19").Value = 0.04487810807300637
$ws.Range("E20").Value = 3
$ws.Range("G20").Value = 296.2543436666667
$ws.Range("H20").Value = 888.763031
$ws.Range("I20").Value = 0.2549833862118423
$ws.Range("J20").Value = 0.2573742821060879
$ws.Range("K20").Value = 3
$ws.Range("M20").Value = 58.81429933333334
$ws.Range("N20").Value = 176.442898
$ws.Range("O20").Value = 0.1539119344629323
$ws.Range("P20").Value = 0.1594685153879524
$ws.Range("Q20").Value = 17423.99164721154
$ws.Range("R20").Value = 156815.9248249038
$ws.Range("S20").Value = 0.03924498622777362
$ws.Range("T20").Value = 0.04104309466649789
$ws.Range("E21").Value = 3
$ws.Range("G21").Value = 296.2543436666667
$ws.Range("H21").Value = 888.763031
$ws.Range("I21").Value = 0.2549833862118423
$ws.Range("J21").Value = 0.2573742821060879
$ws.Range("K21").Value = 2
$ws.Range("M21").Value = 39.945198
$ws.Range("N21").Value = 79.890396
$ws.Range("O21").Value = 0.1045331282761778
$ws.Range("P21").Value = 0.07220467918111168
$ws.Range("Q21").Value = 11833.93841612505
$ws.Range("R21").Value = 71003.63049675027
$ws.Range("S21").Value = 0.0266542110191767
$ws.Range("T21").Value = 0.01858362746893901
$ws.Range("E22").Value = 2
$ws.Range("G22").Value = 32.379461
$ws.Range("H22").Value = 64.758922
$ws.Range("I22").Value = 0.02786870399032479
$ws.Range("J22").Value = 0.0187533464808508
$ws.Range("K22").Value = 3
$ws.Range("M22").Value = 172.9098156666666
$ws.Range("N22").Value = 518.7294469999999
$ws.Range("O22").Value = 0.4524900325013766
$ws.Range("P22").Value = 0.4688259813160829
$ws.Range("Q22").Value = 5598.726632896021
$ws.Range("R22").Value = 33592.35979737613
$ws.Range("S22").Value = 0.01261031077435331
$ws.Range("T22").Value = 0.008792056066845388
$ws.Range("E23").Value = 2
$ws.Range("G23").Value = 32.379461
$ws.Range("H23").Value = 64.758922
$ws.Range("I23").Value = 0.02786870399032479
$ws.Range("J23").Value = 0.0187533464808508
$ws.Range("K23").Value = 3
$ws.Range("M23").Value = 46.15041300000001
$ws.Range("N23").Value = 138.451239
$ws.Range("O23").Value = 0.120771639237527
$ws.Range("P23").Value = 0.1251317779701883
$ws.Range("Q23").Value = 1494.325497867393
$ws.Range("R23").Value = 8965.952987204359
$ws.Range("S23").Value = 0.003365749064336935
$ws.Range("T23").Value = 0.002346639588039835
$ws.Range("E24").Value = 2
$ws.Range("G24").Value = 32.379461
$ws.Range("H24").Value = 64.758922
$ws.Range("I24").Value = 0.02786870399032479
$ws.Range("J24").Value = 0.0187533464808508
$ws.Range("K24").Value = 3
$ws.Range("M24").Value = 64.30983099999999
$ws.Range("N24").Value = 192.929493
$ws.Range("O24").Value = 0.1682932655219863
$ws.Range("P24").Value = 0.1743690461446646
$ws.Range("Q24").Value = 2082.317664781091
$ws.Range("R24").Value = 12493.90598868655
$ws.Range("S24").Value = 0.004690115200397369
$ws.Range("T24").Value = 0.003270003137886358
$ws.Range("E25").Value = 2
$ws.Range("G25").Value = 32.379461
$ws.Range("H25").Value = 64.758922
$ws.Range("I25").Value = 0.02786870399032479
$ws.Range("J25").Value = 0.0187533464808508
$ws.Range("K25").Value = 3
$ws.Range("M25").Value = 58.81429933333334
$ws.Range("N25").Value = 176.442898
$ws.Range("O25").Value = 0.1539119344629323
$ws.Range("P25").Value = 0.1594685153879524
$ws.Range("Q25").Value = 1904.375311505993
$ws.Range("R25").Value = 11426.25186903596
$ws.Range("S25").Value = 0.004289326142125728
$ws.Range("T25").Value = 0.002990568321857159
$ws.Range("E26").Value = 2
$ws.Range("G26").Value = 32.379461
$ws.Range("H26").Value = 64.758922
$ws.Range("I26").Value = 0.02786870399032479
$ws.Range("J26").Value = 0.0187533464808508
$ws.Range("K26").Value = 2
$ws.Range("M26").Value = 39.945198
$ws.Range("N26").Value = 79.890396
$ws.Range("O26").Value = 0.1045331282761778
$ws.Range("P26").Value = 0.07220467918111168
$ws.Range("Q26").Value = 1293.403980778278
$ws.Range("R26").Value = 5173.615923113111
$ws.Range("S26").Value = 0.002913202809111449
$ws.Range("T26").Value = 0.001354079366222062
